# Fixed compute horizontal score algorithm and processing of subheaders.
#
# Row 3 holds the "subheader" labels (SH1, SH3, SH5 ...). The fix removes
# the old "SH3" subheader from D3, and introduces two new subheaders,
# "SH2" (C3) and "SH4" (E3), so the subheader row reads:
#   B3=SH1, C3=SH2, D3=<blank>, E3=SH4, F3=SH5
# Row 4 (the per-column score headers SH21..SH25) is unaffected.
# The active selection moves from D4 to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "SH4"
$ws.Range("C3").Value = "SH2"

# --- Update the active cell selection to E4 ---
$ws.Range("E4").Select()
